$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912340006447"
$wb.Worksheets.Item(2).Name = "NB_TO-1650291236882815"
$wb.Worksheets.Item(3).Name = "RS_TO-1650291236884814"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912369318125"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912369948096"

# --- Sheet 1 (GNG) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912339436536.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912339676428.csv"
$ws1.Range("B4").Value = "go_stims-16502912339686422.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291233999645.csv"

# --- Sheet 2 (NB) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_8-16502912345426433.csv"
$ws2.Range("B3").Value = "ZB-match_9-16502912352878213.csv"
$ws2.Range("B4").Value = "TB-16502912361228454.csv"
$ws2.Range("B5").Value = "OB-1650291235791807.csv"
$ws2.Range("B6").Value = "OB-16502912357438405.csv"
$ws2.Range("B7").Value = "ZB-match_3-1650291234649646.csv"
$ws2.Range("B8").Value = "OB-16502912353758178.csv"
$ws2.Range("B9").Value = "TB-16502912368628085.csv"
$ws2.Range("B10").Value = "TB-16502912364758434.csv"

# --- Sheet 3 (RS) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912368988116.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912368878198.csv"
$ws4.Range("B4").Value = "MM_stims-16502912369148188.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291236899812.csv"
$ws4.Range("B6").Value = "MM_stims-16502912369308164.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912369158094.csv"

# --- Sheet 5 (vSAT) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291236947822.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912369348128.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502912369638133.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912369788134.csv"
